$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.213.70'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '1.860.17'
$ws.Range('E3').Value = '  -1.10%  '
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '0.7028'
$ws.Range('E5').Value = '  -1.95%  '
$ws.Range('D6').Value = '241.62'
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '0.07811'
$ws.Range('E8').Value = '  -1.73%  '
$ws.Range('D9').Value = '0.3109'
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('D10').Value = '23.84'
$ws.Range('E10').Value = '  -4.06%  '
$ws.Range('D11').Value = '0.07802'
$ws.Range('E11').Value = '  -3.44%  '
$ws.Range('D12').Value = '1.843.33'
$ws.Range('E12').Value = '  -2.62%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '5.124'
$ws.Range('E13').Value = '  -1.44%  '
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').Value = '92.55'
$ws.Range('E14').Value = '  -2.17%  '
$ws.Range('D15').Value = '0.6892'
$ws.Range('E15').Value = '  -2.54%  '
$ws.Range('D16').Value = '6.549'
$ws.Range('E16').Value = '  +2.67%  '
$ws.Range('D17').Value = '0.000008440'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').Value = '29.203.34'
$ws.Range('E18').Value = '  -1.61%  '
$ws.Range('D19').Value = '249.96'
$ws.Range('D20').Value = '2.100.97'
$ws.Range('E20').Value = '  -0.94%  '
$ws.Range('D21').Value = '12.91'
$ws.Range('E21').Value = '  -3.10%  '
$ws.Range('D22').Value = '1.0000'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').Value = '7.590'
$ws.Range('E23').Value = '  -0.91%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('E25').Value = '  -2.71%  '
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('D27').Value = '8.888'
$ws.Range('E27').Value = '  -1.83%  '
$ws.Range('D28').Value = '18.56'
$ws.Range('E28').Value = '  -1.96%  '
$ws.Range('E29').Value = '  +3.80%  '
$ws.Range('D30').Value = '4.277'
$ws.Range('E30').Value = '  -3.00%  '
$ws.Range('D31').Value = '4.253'
$ws.Range('E31').Value = '  -1.36%  '
$ws.Range('D32').Value = '1.206'
$ws.Range('E32').Value = '  -2.42%  '
$ws.Range('D33').Value = '0.05212'
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('D34').Value = '0.7589'
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('D35').Value = '1.868'
$ws.Range('E35').Value = '  -3.32%  '
$ws.Range('D36').Value = '1.175'
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').Value = '2.709'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').Value = '0.01861'
$ws.Range('E38').Value = '  -1.05%  '
$ws.Range('D39').Value = '1.226.84'
$ws.Range('E39').Value = '  -4.73%  '
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('D41').Value = '0.8995'
$ws.Range('E41').Value = '  -0.74%  '
$ws.Range('D42').Value = '109.61'
$ws.Range('E42').Value = '  -1.62%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = '0.9996'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '5.728'
$ws.Range('E44').Value = '  -10.60%  '
$ws.Range('D45').Value = '2.001.39'
$ws.Range('E45').Value = '  -1.48%  '
$ws.Range('E46').Value = '  -4.93%  '
$ws.Range('D47').Value = '65.50'
$ws.Range('E47').Value = '  -11.34%  '
$ws.Range('D48').Value = '0.5187'
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('D49').Value = '9.526'
$ws.Range('E49').Value = '  +0.30%  '
$ws.Range('E50').Value = '  -2.47%  '
$ws.Range('D51').Value = '7.035'
$ws.Range('E51').Value = '  -0.65%  '
